# Add a new "PROJECT LINK" slide (with the GitHub repo link) right before
# the final "RESULTS" slide. This mirrors inserting a new slide at
# position 12 (pushing the existing last slide to position 13) using the
# "Title and Content" layout (the same slide layout used by RESULTS).

$p = $ppt.ActivePresentation

# ppLayoutText (2) == the "Title and Content" slideLayout2.xml used by
# the neighbouring RESULTS slide.
$s = $p.Slides.Add(12, 2)

$title = $s.Shapes.Item(1)
$body  = $s.Shapes.Item(2)

$title.Name = "Title 1"
$body.Name  = "Text Placeholder 2"

# --- Title -----------------------------------------------------------
$title.TextFrame.TextRange.Text = "PROJECT LINK"

# --- Body placeholder --------------------------------------------------
# Resize/reposition to match the target layout (values are in points;
# EMU = points * 12700).
$body.Left   = 48.00003937007874
$body.Top    = 124.20003937007874
$body.Width  = 864.0000393700788
$body.Height = 109.05476377952756

$bodyTr = $body.TextFrame.TextRange
$bodyTr.Text = "https://github.com/Saijyothula9/keylogger_project.git"
# Prepend four blank paragraphs ahead of the link paragraph.
$bodyTr.InsertBefore("`r`r`r`r") | Out-Null
# Centre every paragraph in the placeholder.
$bodyTr.ParagraphFormat.Alignment = 2

# Re-fetch the (now five-paragraph) range and hyperlink the last one.
$fresh = $body.TextFrame.TextRange
$linkPara = $fresh.Paragraphs(5, 1)
$linkPara.ActionSettings.Item(1).Hyperlink.Address = "https://github.com/Saijyothula9/keylogger_project.git"

Write-Host "Inserted PROJECT LINK slide at index 12 (slide count now $($p.Slides.Count))."
